# Update Name of Algo
# Apply updated RandomForest imputation results to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -11.5184
$ws.Range("C3").Value = -12.0697
$ws.Range("C5").Value = -13.059
$ws.Range("E7").Value = 11.9977
$ws.Range("B9").Value = 8.740300000000005
$ws.Range("E9").Value = 14.52300000000001
$ws.Range("C11").Value = -13.65010000000001
$ws.Range("C12").Value = -14.52240000000002
$ws.Range("B13").Value = 6.603300000000004
$ws.Range("B16").Value = 8.978700000000009
$ws.Range("B18").Value = 6.198300000000002
$ws.Range("B20").Value = 5.389499999999998
$ws.Range("C21").Value = -13.26360000000001
$ws.Range("E21").Value = 12.74709999999998
